$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "porta" (port) column E annotations, entered in numeric port order 0-7 ---
# Rows 9:10 share port 1, rows 14:15 share port 2, rows 24:25 share port 7 -> merged cells.
# Single-row port annotations (not merged) go in rows 27-31.
$ws.Range("E27").Value = "porta 0"

$ws.Range("E9").Value = "porta 1"
$ws.Range("E9:E10").Merge()

$ws.Range("E14").Value = "porta 2"
$ws.Range("E14:E15").Merge()

$ws.Range("E28").Value = "porta 3"
$ws.Range("E29").Value = "porta 4"
$ws.Range("E30").Value = "porta 5"
$ws.Range("E31").Value = "porta 6"

$ws.Range("E24").Value = "porta 7"
$ws.Range("E24:E25").Merge()

# --- Rename the four pin labels that lost their "O" (IOxx -> Ixx) for pins 4-7 ---
$ws.Range("B5").Value = "I36"
$ws.Range("B6").Value = "I39"
$ws.Range("B7").Value = "I34"
$ws.Range("B8").Value = "I35"

# Format the new E column cells: centered, Century Gothic 12pt, no fill (matches the
# rest of the table's font/alignment but without the colored banding used in B:D)
$eRanges = @("E9:E10", "E14:E15", "E24:E25", "E27", "E28", "E29", "E30", "E31")
foreach ($addr in $eRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Century Gothic"
    $rng.Font.Size = 12
    $rng.HorizontalAlignment = -4108
    $rng.Interior.Pattern = -4142
}

# --- View state: scroll back to the top and move the active selection ---
$ws.Activate()
$ws.Range("J16").Select()
